$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new row before row 20. This shifts the (currently empty)
#    rows 20-23 and the signature block (rows 24,25) down by one, to
#    rows 21-24 and 25,26 respectively, making room for a new data row.
# ------------------------------------------------------------------
$ws.Rows("20:20").Insert() | Out-Null

# ------------------------------------------------------------------
# 2) Fix up formatting: row 20 becomes the new "last row" of the table
#    (bottom border), so copy the format that row 19 currently has.
#    Row 19 becomes a normal "middle" row, so copy the format from
#    row 18 onto it (must happen AFTER copying row19's original format
#    down to row 20).
# ------------------------------------------------------------------
$ws.Range("B19:J19").Copy() | Out-Null
$ws.Range("B20:J20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B18:J18").Copy() | Out-Null
$ws.Range("B19:J19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Re-sort the "Periodo Mora" column from descending to ascending,
#    and add the new period 2508 as the new last row.
# ------------------------------------------------------------------
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"
$ws.Range("E19").Value = "2507"
$ws.Range("E20").Value = "2508"

# ------------------------------------------------------------------
# 4) Fill in the rest of the new row 20 (same worker/employer data as
#    the other rows, same Salario Basico / Valor Mora amounts).
# ------------------------------------------------------------------
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047399469"
$ws.Range("D20").Value = "MARIA DEL ROSARIO CARRILLO ANGULO"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# ------------------------------------------------------------------
# 5) Update the summary figures: total Valor Mora and Cant. Periodos.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 284700
$ws.Range("F13").Value = 5

Write-Host "done"
